$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 (5702a499-...md) status moves from
#     "Ready for handoff" to "Handed back: in sync with en-US" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row 3 handback recorded ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K3").Value = "2016-09-04 10:53:34"
$wsZhCn.Range("P3").Value = ""

# --- de-de sheet: row 3 handback recorded ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K3").Value = "2016-09-04 10:53:41"
$wsDeDe.Range("P3").Value = ""

# Error Detail column narrows now that the long error text is gone
$wsZhCn.Columns("P").ColumnWidth = 12.8
$wsDeDe.Columns("P").ColumnWidth = 12.8
